$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily data rows (2021-09-02 through 2021-09-09 / serials 44441-44448)
$data = @(
    @(44441, 4, 10, 623.4413965087282),
    @(44442, 1, 11, 685.785536159601),
    @(44443, 0, 6, 374.0648379052369),
    @(44444, 0, 6, 374.0648379052369),
    @(44445, 0, 6, 374.0648379052369),
    @(44446, 0, 5, 311.7206982543641),
    @(44447, 0, 5, 311.7206982543641),
    @(44448, 0, 1, 62.34413965087282)
)

$startRow = 367
$lastRow = $startRow - 1

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]

    # Copy the formatting (style) of the last existing date cell (column A)
    # down onto the new date cell, mirroring how the sheet keeps extending.
    $ws.Range("A$lastRow").Copy()
    $ws.Range("A$row").PasteSpecial(-4122)

    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
}

$excel.CutCopyMode = 0
